# Add two new slides ("Differing Slide" and "Second slide"), both using the
# "Title and Content" layout (layout index 2 == ppLayoutText's custom-layout
# equivalent on this deck, matches slideLayout2.xml).
#
# Both Slides.Add calls target index 2 (right after the existing slide 1), so
# the slide added *second* ends up in final position 2, and the slide added
# *first* is pushed down to final position 3 - this reproduces the out-of-
# sequence slide ids (256, 258, 257) seen in the target deck.

$p = $ppt.ActivePresentation

# Added first -> ends up as the final slide (position 3, "Differing Slide").
$sDiffering = $p.Slides.Add(2, 2)
$sDiffering.Shapes.Item(1).TextFrame.TextRange.Text = "Differing Slide`t"
$diffContent = $sDiffering.Shapes.Item(2).TextFrame.TextRange
$diffContent.Text = "Differing "
$diffContent.InsertAfter("contnet") | Out-Null

# Added second -> ends up right after slide 1 (position 2, "Second slide").
$sSecond = $p.Slides.Add(2, 2)
$sSecond.Shapes.Item(1).TextFrame.TextRange.Text = "Second slide"
$secondContent = $sSecond.Shapes.Item(2).TextFrame.TextRange
$secondContent.Text = "This is text on "
$secondContent.InsertAfter("second slide") | Out-Null

Write-Output ("Slides.Count=" + $p.Slides.Count)
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    Write-Output ("Slide " + $i + " id=" + $sl.SlideID + " title=" + $sl.Shapes.Item(1).TextFrame.TextRange.Text)
}
